$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New expert result: KI (Nov-11-2023) -- adds three rows (14-16) to the tracker.
# Columns: A=Task_type, B=Date, C=expert, D=submission_file_name,
#          E=response_collected, F=Further_process

# The "Date" column holds values like "Nov-11-2023" which Excel's smart
# text-to-date recognition would otherwise silently convert into a date
# serial number. Force those three cells to Text first so the literal
# string is stored (matching every other Date cell in this sheet), then
# drop the now-unneeded explicit number format so the cells keep the
# sheet's default (General) styling.
$ws.Range("B14:B16").NumberFormat = "@"

# Row 14
$ws.Range("A14").Value = "pairwise"
$ws.Range("B14").Value = "Nov-11-2023"
$ws.Range("C14").Value = "KI"
$ws.Range("D14").Value = "all_submitted_tracker_KI_Nov-11-2023.csv"
$ws.Range("E14").Value = "master_all_responses_KINov-11-2023.csv"

# Row 15
$ws.Range("A15").Value = "paiewise_resub"
$ws.Range("B15").Value = "Nov-11-2023"
$ws.Range("C15").Value = "KI"
$ws.Range("D15").Value = "re_submitted_tracker_KINov-11-2023.csv"
$ws.Range("E15").Value = "master_all_responses_KI_resubNov-11-2023.csv"

# Row 16
$ws.Range("A16").Value = "paiewise_resub"
$ws.Range("B16").Value = "Nov-11-2023"
$ws.Range("C16").Value = "KI"
$ws.Range("D16").Value = "re_submitted_tracker_KI2Nov-11-2023.csv"
$ws.Range("E16").Value = "master_all_responses_KI_resub2Nov-11-2023.csv"
$ws.Range("F16").Value = "master_all_responses_KINov-11-2023_to_KI_resub2Nov-11-2023_Kiyomi.csv"

$ws.Range("B14:B16").ClearFormats()

$ws.Range("E19").Select()

$wb.Save()
